$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.893.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.545.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.27'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.46%  '
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0856'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.766.67'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.550.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.864.31'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.45%  '
$ws.Range("E24").Value = '  -3.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.361.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.964'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.78%  '
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0164'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.518'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.804'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.987'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.73'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.681.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0504'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("E50").Value = '  -0.36%  '
